# Update the public EPEX Spot prices workbook with the latest day of data
# (03-aug on "Prix Spot", 2025-08-01 on "Gaz" and "CO2").

$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column AY (03-aug) -----------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the header formatting from the previous day's header cell (AX1) so the
# new header cell (AY1) keeps the same bold/centered/bordered style, then set
# its text.
$wsPrix.Range("AX1").Copy()
$wsPrix.Range("AY1").PasteSpecial(-4122)
$wsPrix.Range("AY1").Value = "03-aug"

$prixValues = @(
    80.40000000000001,
    69.26000000000001,
    42.96,
    37.38,
    30.02,
    30.2,
    32.03,
    26.28,
    6.16,
    0,
    -1.01,
    -2.1,
    -3,
    -9.9,
    -10.08,
    -4.98,
    -1.49,
    -0.01,
    12.51,
    40,
    78.2,
    96.13,
    97.40000000000001,
    85.09
)

for ($i = 0; $i -lt $prixValues.Count; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 51).Value = $prixValues[$i]
}

# --- Sheet "Gaz": append row 48 (2025-08-01) --------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date cell to stay plain text ("2025-08-01") instead of being
# auto-converted to a date serial number, then drop back to the default
# (unstyled) cell style so it matches the other rows in the column.
$wsGaz.Range("A48").NumberFormat = "@"
$wsGaz.Range("A48").Value = "2025-08-01"
$wsGaz.Range("A48").Style = "Normal"

$wsGaz.Range("B48").Value = 32.65

# --- Sheet "CO2": append row 48 (2025-08-01) --------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A48").NumberFormat = "@"
$wsCo2.Range("A48").Value = "2025-08-01"
$wsCo2.Range("A48").Style = "Normal"

$wsCo2.Range("B48").Value = 70.58
